# Auto-generated edit script applying the Anima_Profits leve-price refresh
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets), per scheduled-runner diff.
$wb = $excel.ActiveWorkbook

function Set-LeveCell {
    param($ws, [int]$row, [int]$col, $value)
    if ($null -eq $value) {
        $ws.Cells.Item($row, $col).Value = $null
    } else {
        $ws.Cells.Item($row, $col).Value = $value
    }
}

$ws = $wb.Worksheets.Item("ALC")
Set-LeveCell $ws 42 8 331.33334
Set-LeveCell $ws 42 10 402
Set-LeveCell $ws 42 12 1206
Set-LeveCell $ws 42 14 -1666

$ws = $wb.Worksheets.Item("ALC")
Set-LeveCell $ws 69 8 2854.1428
Set-LeveCell $ws 69 9 2660
Set-LeveCell $ws 69 10 2999.75
Set-LeveCell $ws 69 11 7980
Set-LeveCell $ws 69 12 8999.25
Set-LeveCell $ws 69 13 -7106
Set-LeveCell $ws 69 14 -10747.25

$ws = $wb.Worksheets.Item("ALC")
Set-LeveCell $ws 72 8 2854.1428
Set-LeveCell $ws 72 9 2660
Set-LeveCell $ws 72 10 2999.75
Set-LeveCell $ws 72 11 23940
Set-LeveCell $ws 72 12 26997.75
Set-LeveCell $ws 72 13 -19572
Set-LeveCell $ws 72 14 -35733.75

$ws = $wb.Worksheets.Item("ALC")
Set-LeveCell $ws 137 8 1293.8846
Set-LeveCell $ws 137 9 1149.2941
Set-LeveCell $ws 137 10 1567
Set-LeveCell $ws 137 11 3447.8823
Set-LeveCell $ws 137 12 4701
Set-LeveCell $ws 137 13 -897.8823000000002
Set-LeveCell $ws 137 14 -9801

$ws = $wb.Worksheets.Item("ALC")
Set-LeveCell $ws 138 8 2617.76
Set-LeveCell $ws 138 9 2841.7144
Set-LeveCell $ws 138 10 2530.6667
Set-LeveCell $ws 138 11 8525.143199999999
Set-LeveCell $ws 138 12 7592.000100000001
Set-LeveCell $ws 138 13 -3385.143199999999
Set-LeveCell $ws 138 14 -17872.0001

$ws = $wb.Worksheets.Item("ARM")
Set-LeveCell $ws 54 8 0
Set-LeveCell $ws 54 10 0
Set-LeveCell $ws 54 14 $null

$ws = $wb.Worksheets.Item("ARM")
Set-LeveCell $ws 132 8 4307.1113
Set-LeveCell $ws 132 9 4616.926
Set-LeveCell $ws 132 10 3377.6667
Set-LeveCell $ws 132 11 13850.778
Set-LeveCell $ws 132 12 10133.0001
Set-LeveCell $ws 132 13 -11320.778
Set-LeveCell $ws 132 14 -15193.0001

$ws = $wb.Worksheets.Item("BSM")
Set-LeveCell $ws 128 8 3485
Set-LeveCell $ws 128 9 3485
Set-LeveCell $ws 128 11 10455
Set-LeveCell $ws 128 13 -7965

$ws = $wb.Worksheets.Item("BSM")
Set-LeveCell $ws 134 8 2864.4375
Set-LeveCell $ws 134 9 2710.6365
Set-LeveCell $ws 134 10 3202.8
Set-LeveCell $ws 134 11 8131.9095
Set-LeveCell $ws 134 12 9608.400000000001
Set-LeveCell $ws 134 13 -5596.9095
Set-LeveCell $ws 134 14 -14678.4

$ws = $wb.Worksheets.Item("CRP")
Set-LeveCell $ws 31 8 4356.595
Set-LeveCell $ws 31 9 1294.1428
Set-LeveCell $ws 31 10 7419.048
Set-LeveCell $ws 31 11 1294.1428
Set-LeveCell $ws 31 12 7419.048
Set-LeveCell $ws 31 13 -999.1428000000001
Set-LeveCell $ws 31 14 -8009.048

$ws = $wb.Worksheets.Item("CRP")
Set-LeveCell $ws 34 8 4356.595
Set-LeveCell $ws 34 9 1294.1428
Set-LeveCell $ws 34 10 7419.048
Set-LeveCell $ws 34 11 1294.1428
Set-LeveCell $ws 34 12 7419.048
Set-LeveCell $ws 34 13 -1092.1428
Set-LeveCell $ws 34 14 -7823.048

$ws = $wb.Worksheets.Item("CRP")
Set-LeveCell $ws 58 8 2570.9285
Set-LeveCell $ws 58 9 2349.3
Set-LeveCell $ws 58 10 3125
Set-LeveCell $ws 58 11 2349.3
Set-LeveCell $ws 58 12 3125
Set-LeveCell $ws 58 13 -2146.3
Set-LeveCell $ws 58 14 -3531

$ws = $wb.Worksheets.Item("CRP")
Set-LeveCell $ws 132 8 11907133
Set-LeveCell $ws 132 9 1249.5714
Set-LeveCell $ws 132 10 23813018
Set-LeveCell $ws 132 11 3748.7142
Set-LeveCell $ws 132 12 71439054
Set-LeveCell $ws 132 13 -1218.7142
Set-LeveCell $ws 132 14 -71444114

$ws = $wb.Worksheets.Item("CRP")
Set-LeveCell $ws 134 8 1100.591
Set-LeveCell $ws 134 9 811.05554
Set-LeveCell $ws 134 10 2403.5
Set-LeveCell $ws 134 11 2433.16662
Set-LeveCell $ws 134 12 7210.5
Set-LeveCell $ws 134 13 101.83338
Set-LeveCell $ws 134 14 -12280.5

$ws = $wb.Worksheets.Item("CRP")
Set-LeveCell $ws 136 8 2570.9285
Set-LeveCell $ws 136 9 2349.3
Set-LeveCell $ws 136 10 3125
Set-LeveCell $ws 136 11 7047.900000000001
Set-LeveCell $ws 136 12 9375
Set-LeveCell $ws 136 13 -4497.900000000001
Set-LeveCell $ws 136 14 -14475

$ws = $wb.Worksheets.Item("CUL")
Set-LeveCell $ws 5 8 508.57895
Set-LeveCell $ws 5 9 405.9655
Set-LeveCell $ws 5 10 614.8570999999999
Set-LeveCell $ws 5 11 1217.8965
Set-LeveCell $ws 5 12 1844.5713
Set-LeveCell $ws 5 13 -1105.8965
Set-LeveCell $ws 5 14 -2068.5713

$ws = $wb.Worksheets.Item("CUL")
Set-LeveCell $ws 38 8 1495.5555
Set-LeveCell $ws 38 9 80
Set-LeveCell $ws 38 10 1900
Set-LeveCell $ws 38 11 240
Set-LeveCell $ws 38 12 5700
Set-LeveCell $ws 38 13 107
Set-LeveCell $ws 38 14 -6394

$ws = $wb.Worksheets.Item("CUL")
Set-LeveCell $ws 68 8 2329.4194
Set-LeveCell $ws 68 9 3153.054
Set-LeveCell $ws 68 10 1785.2322
Set-LeveCell $ws 68 11 9459.162
Set-LeveCell $ws 68 12 5355.696599999999
Set-LeveCell $ws 68 13 -8648.162
Set-LeveCell $ws 68 14 -6977.696599999999

$ws = $wb.Worksheets.Item("CUL")
Set-LeveCell $ws 71 8 2329.4194
Set-LeveCell $ws 71 9 3153.054
Set-LeveCell $ws 71 10 1785.2322
Set-LeveCell $ws 71 11 28377.486
Set-LeveCell $ws 71 12 16067.0898
Set-LeveCell $ws 71 13 -24321.486
Set-LeveCell $ws 71 14 -24179.0898

$ws = $wb.Worksheets.Item("CUL")
Set-LeveCell $ws 98 8 333526.66
Set-LeveCell $ws 98 9 290
Set-LeveCell $ws 98 10 1000000
Set-LeveCell $ws 98 11 870
Set-LeveCell $ws 98 12 3000000
Set-LeveCell $ws 98 13 628
Set-LeveCell $ws 98 14 -3002996

$ws = $wb.Worksheets.Item("CUL")
Set-LeveCell $ws 107 8 1616.2059
Set-LeveCell $ws 107 9 329.91666
Set-LeveCell $ws 107 10 2317.818
Set-LeveCell $ws 107 11 989.7499799999999
Set-LeveCell $ws 107 12 6953.454000000001
Set-LeveCell $ws 107 13 930.2500200000001
Set-LeveCell $ws 107 14 -10793.454

$ws = $wb.Worksheets.Item("CUL")
Set-LeveCell $ws 135 8 508.57895
Set-LeveCell $ws 135 9 405.9655
Set-LeveCell $ws 135 10 614.8570999999999
Set-LeveCell $ws 135 11 3653.6895
Set-LeveCell $ws 135 12 5533.7139
Set-LeveCell $ws 135 13 -1118.6895
Set-LeveCell $ws 135 14 -10603.7139

$ws = $wb.Worksheets.Item("GSM")
Set-LeveCell $ws 41 8 3124
Set-LeveCell $ws 41 9 3124
Set-LeveCell $ws 41 11 3124
Set-LeveCell $ws 41 13 -2769

$ws = $wb.Worksheets.Item("GSM")
Set-LeveCell $ws 132 8 2899
Set-LeveCell $ws 132 9 2947.2
Set-LeveCell $ws 132 10 2872.2222
Set-LeveCell $ws 132 11 8841.599999999999
Set-LeveCell $ws 132 12 8616.6666
Set-LeveCell $ws 132 13 -6311.599999999999
Set-LeveCell $ws 132 14 -13676.6666

$ws = $wb.Worksheets.Item("LTW")
Set-LeveCell $ws 132 8 2828.1628
Set-LeveCell $ws 132 9 2419.1785
Set-LeveCell $ws 132 10 3591.6
Set-LeveCell $ws 132 11 7257.5355
Set-LeveCell $ws 132 12 10774.8
Set-LeveCell $ws 132 13 -4727.5355
Set-LeveCell $ws 132 14 -15834.8

$ws = $wb.Worksheets.Item("WVR")
Set-LeveCell $ws 3 8 25230.285
Set-LeveCell $ws 3 9 6000
Set-LeveCell $ws 3 10 32922.4
Set-LeveCell $ws 3 11 6000
Set-LeveCell $ws 3 12 32922.4
Set-LeveCell $ws 3 13 -5886
Set-LeveCell $ws 3 14 -33150.4

$ws = $wb.Worksheets.Item("WVR")
Set-LeveCell $ws 81 8 7067.143
Set-LeveCell $ws 81 9 8534
Set-LeveCell $ws 81 10 3400
Set-LeveCell $ws 81 11 17068
Set-LeveCell $ws 81 12 6800
Set-LeveCell $ws 81 13 -16007
Set-LeveCell $ws 81 14 -8922

$ws = $wb.Worksheets.Item("WVR")
Set-LeveCell $ws 84 8 7067.143
Set-LeveCell $ws 84 9 8534
Set-LeveCell $ws 84 10 3400
Set-LeveCell $ws 84 11 85340
Set-LeveCell $ws 84 12 34000
Set-LeveCell $ws 84 13 -80036
Set-LeveCell $ws 84 14 -44608

$ws = $wb.Worksheets.Item("WVR")
Set-LeveCell $ws 129 8 100429
Set-LeveCell $ws 129 10 100429
Set-LeveCell $ws 129 12 100429
Set-LeveCell $ws 129 14 -110429

$ws = $wb.Worksheets.Item("WVR")
Set-LeveCell $ws 132 8 5611395.5
Set-LeveCell $ws 132 9 2805.2593
Set-LeveCell $ws 132 10 11668673
Set-LeveCell $ws 132 11 8415.777900000001
Set-LeveCell $ws 132 12 35006019
Set-LeveCell $ws 132 13 -5885.777900000001
Set-LeveCell $ws 132 14 -35011079

$ws = $wb.Worksheets.Item("WVR")
Set-LeveCell $ws 136 8 2464.0876
Set-LeveCell $ws 136 9 2291.561
Set-LeveCell $ws 136 10 2906.1875
Set-LeveCell $ws 136 11 6874.683000000001
Set-LeveCell $ws 136 12 8718.5625
Set-LeveCell $ws 136 13 -4324.683000000001
Set-LeveCell $ws 136 14 -13818.5625
